$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison -> column D (MyForecast) updates ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("D3").Value = 102
$wsForecast.Range("D4").Value = 99
$wsForecast.Range("D5").Value = 109
$wsForecast.Range("D9").Value = 97
$wsForecast.Range("D10").Value = 88
$wsForecast.Range("D11").Value = 96
$wsForecast.Range("D14").Value = 104
$wsForecast.Range("D15").Value = 93
$wsForecast.Range("D16").Value = 86
$wsForecast.Range("D17").Value = 88

# --- Sheet: Summary -> column B updates (values are stored as text) ---
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "1612"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "846"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "415"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "86"
